$d = $word.ActiveDocument

$pairs = @(
    @("56×27=", "34×39="),
    @("55×59=", "47×60="),
    @("60×18=", "96×99="),
    @("63×44=", "82×91="),
    @("37×61=", "97×87="),
    @("80×13=", "51×67="),
    @("36×99=", "63×62="),
    @("13×20=", "46×44="),
    @("42×60=", "61×76="),
    @("23×27=", "23×44="),
    @("23×22=", "49×67="),
    @("22×59=", "71×26="),
    @("55×33=", "51×40="),
    @("74×53=", "59×90="),
    @("30×99=", "38×51="),
    @("77×58=", "72×72="),
    @("83×16=", "62×36="),
    @("17×90=", "21×12="),
    @("27×42=", "17×33="),
    @("25×53=", "80×12="),
    @("45×55=", "54×96="),
    @("76×69=", "80×45="),
    @("66×41=", "90×22="),
    @("43×68=", "38×62="),
    @("44×88=", "23×24=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
